# "updated the experiment template"
# Fill in the answer column (B) with a "-" placeholder for every
# question row (rows 2-7; B1 already had its answer), then move the
# selection to the first empty cell below the filled data (B8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B7").Value = "-"

$null = $ws.Range("B8").Select()
